$d = $word.ActiveDocument

# --- Change 1: remove the "Meta description" paragraph that used to follow the
#     Heading1 title paragraph at the very top of the document. ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- Change 2: at the end of the document, insert a new bold paragraph
#     ("Play East Sea Dragon King for Free | Review and Gameplay") right
#     before the final (italic) paragraph, and replace that final paragraph's
#     text (was the "Create a feature image..." image prompt) with the old
#     meta-description copy, keeping its italic formatting. ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$insertStart = $lastPara.Range.Start
$insertPoint = $d.Range($insertStart, $insertStart)

$newHeadingText = "Play East Sea Dragon King for Free | Review and Gameplay"
$xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $newHeadingText + '</w:t></w:r></w:p>'
$insertPoint.InsertXML($xmlFragment)

$breakPos = $insertStart + $newHeadingText.Length
$breakRange = $d.Range($breakPos, $breakPos)
$breakRange.InsertParagraphAfter()

$oldPrompt = 'Create a feature image for the game ' + [char]0x22 + 'East Sea Dragon King' + [char]0x22 + ' that captures the Asian and cartoon theme of the game. The image should prominently feature a happy Maya warrior with glasses, fitting in with the overall aesthetic of the game. The image should be bright and eye-catching, with a colour palette that reflects the underwater setting and elements of traditional Asian design. The overall style should be cartoonish and fun, with a strong emphasis on the character of the Maya warrior.'
$newCopy = 'Get a chance to win up to 6,000 times your bet with East Sea Dragon King. Play for free and learn about its features in our comprehensive review.'

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newCopy, 2)
